# Update repository-name formatting from "owner_repo" (underscore) style to
# "owner/repo" (slash) style on both worksheets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("PyTreesRos_model_analysis")
$ws2 = $wb.Worksheets.Item("aggregated_data")

# --- Sheet "PyTreesRos_model_analysis" : column B (repository), rows 2-21 ---
$sheet1B = @{
  2  = "braineniac/robotics-player"
  3  = "jotix16/Robotics-Behaviour-Planning"
  4  = "jotix16/Robotics-Behaviour-Planning"
  5  = "jotix16/Robotics-Behaviour-Planning"
  6  = "KKalem/sam_march"
  7  = "peterheim1/gizmo"
  8  = "peterheim1/gizmo"
  9  = "peterheim1/gizmo"
  10 = "peterheim1/gizmo"
  11 = "peterheim1/gizmo"
  12 = "peterheim1/gizmo"
  13 = "peterheim1/gizmo"
  14 = "peterheim1/gizmo"
  15 = "efills-project/refills_second_review"
  16 = "samiamlabs/dyno"
  17 = "samiamlabs/dyno"
  18 = "simutisernestas/mobile_robot_project"
  19 = "smarc-project/smarc_missions"
  20 = "smarc-project/smarc_missions"
  21 = "Taospirit/roborts_project"
}

foreach ($row in $sheet1B.Keys) {
    $ws1.Cells.Item($row, 2).Value = $sheet1B[$row]
}

# --- Sheet "aggregated_data" : column A (repository), rows 2-10 ---
$sheet2A = @{
  2  = "Kkalem/sam_march"
  3  = "smarc-project/smarc_missions"
  4  = "simutisernestas/mobile_robot_project"
  5  = "samiamlabs/dyno"
  6  = "braineniac/robotics-player"
  7  = "peterheim1/gizmo"
  8  = "Taospirit/roborts_project"
  9  = "jotix16/Robotics-Behaviour-Planning"
  10 = "refills-projectrefills_second_review"
}

foreach ($row in $sheet2A.Keys) {
    $ws2.Cells.Item($row, 1).Value = $sheet2A[$row]
}
